# Apply updated metric values to columns B (currentcaaf_strategy) and
# H (maxdiv_strategy), plus a few C/I cells, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 2164011.04
$ws.Range("H4").Value = 3436474.83
$ws.Range("B5").Value = 6.95
$ws.Range("H5").Value = 7.28
$ws.Range("H8").Value = 0.83
$ws.Range("I8").Value = 0.7
$ws.Range("B9").Value = 1.55
$ws.Range("H9").Value = 1.61
$ws.Range("C10").Value = 1.16
$ws.Range("H10").Value = 1.43
$ws.Range("I10").Value = 1.17
$ws.Range("H11").Value = 1.14
$ws.Range("H12").Value = 1.01
$ws.Range("I12").Value = 0.82
$ws.Range("B13").Value = 2.06
$ws.Range("C13").Value = 2.06
$ws.Range("H13").Value = 2.11
$ws.Range("I13").Value = 2.11
$ws.Range("B14").Value = -47.83
$ws.Range("H14").Value = -44.84
$ws.Range("B15").Value = 2222
$ws.Range("H15").Value = 2130
$ws.Range("B16").Value = 7.77
$ws.Range("H16").Value = 7.87
$ws.Range("B17").Value = 0.86
$ws.Range("C17").Value = 0.86
$ws.Range("H17").Value = 0.6899999999999999
$ws.Range("I17").Value = 0.6899999999999999
$ws.Range("H18").Value = -0
$ws.Range("I18").Value = -0
$ws.Range("B20").Value = 0.22
$ws.Range("H20").Value = 0.25
$ws.Range("B21").Value = 7.77
$ws.Range("H21").Value = 7.41
$ws.Range("B22").Value = 0.5600000000000001
$ws.Range("H22").Value = 0.59
$ws.Range("B23").Value = 0.5600000000000001
$ws.Range("H23").Value = 0.59
$ws.Range("B24").Value = 6.93
$ws.Range("H24").Value = 7.26
$ws.Range("B25").Value = 31.22
$ws.Range("C25").Value = 27.3
$ws.Range("H25").Value = 31.74
$ws.Range("I25").Value = 27.63
$ws.Range("B27").Value = -3.1
$ws.Range("H27").Value = -3.12
$ws.Range("B28").Value = -3.1
$ws.Range("H28").Value = -3.12
$ws.Range("H29").Value = 15
$ws.Range("B31").Value = 1.06
$ws.Range("H31").Value = 1.11
$ws.Range("B32").Value = 1.06
$ws.Range("H32").Value = 1.11
$ws.Range("B33").Value = 1.18
$ws.Range("C33").Value = 1.06
$ws.Range("H33").Value = 1.2
$ws.Range("I33").Value = 1.08
$ws.Range("B34").Value = 2.06
$ws.Range("H34").Value = 2.11
$ws.Range("B35").Value = 2.83
$ws.Range("H35").Value = 2.86
$ws.Range("B36").Value = 1.52
$ws.Range("C36").Value = 1.25
$ws.Range("H36").Value = 1.59
$ws.Range("I36").Value = 1.27
$ws.Range("B37").Value = 1.38
$ws.Range("H37").Value = 1.36
$ws.Range("B38").Value = 3.51
$ws.Range("C38").Value = 3.02
$ws.Range("H38").Value = 3.46
$ws.Range("I38").Value = 3.04
$ws.Range("B39").Value = 3.86
$ws.Range("H39").Value = 4.01
$ws.Range("I39").Value = 3.19
$ws.Range("B41").Value = 2.05
$ws.Range("B42").Value = 3.28
$ws.Range("B43").Value = 7.58
$ws.Range("B44").Value = 1.63
$ws.Range("H44").Value = 3.26
$ws.Range("B45").Value = 1.21
$ws.Range("H45").Value = 1.21
$ws.Range("B46").Value = 5.31
$ws.Range("H46").Value = 6.31
$ws.Range("B47").Value = 4.99
$ws.Range("H47").Value = 5.42
$ws.Range("B48").Value = 6.95
$ws.Range("H48").Value = 7.28
$ws.Range("B49").Value = 18.31
$ws.Range("H49").Value = 19.19
$ws.Range("H50").Value = -12.54
$ws.Range("B51").Value = 18.31
$ws.Range("H51").Value = 19.19
$ws.Range("H52").Value = -12.54
$ws.Range("B53").Value = 35.63
$ws.Range("H53").Value = 37.54
$ws.Range("B54").Value = -21.63
$ws.Range("H54").Value = -20.68
$ws.Range("B55").Value = -2.99
$ws.Range("B56").Value = 155
$ws.Range("H56").Value = 167
$ws.Range("B57").Value = 45243.83
$ws.Range("H57").Value = 76644.05
$ws.Range("H58").Value = 0.06
$ws.Range("B59").Value = 70931.17999999999
$ws.Range("H59").Value = 128296.44
$ws.Range("C60").Value = 2.23
$ws.Range("H60").Value = 1.99
$ws.Range("I60").Value = 2.25
$ws.Range("B61").Value = -1.64
$ws.Range("C61").Value = -2.09
$ws.Range("H61").Value = -1.66
$ws.Range("I61").Value = -2.09
$ws.Range("B62").Value = 62.79
$ws.Range("H62").Value = 62.79
$ws.Range("B63").Value = 62.79
$ws.Range("H63").Value = 62.79
$ws.Range("B64").Value = 69.92
$ws.Range("H64").Value = 70.42
$ws.Range("B65").Value = 79.19
$ws.Range("H65").Value = 81.20999999999999
$ws.Range("B66").Value = 0.76
$ws.Range("H66").Value = 0.6899999999999999
$ws.Range("H67").Value = 0.02
$ws.Range("B68").Value = 92.78
$ws.Range("H68").Value = 83
$ws.Range("B69").Value = 2843176.11
$ws.Range("H69").Value = 4980371.15
